$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from G1 (the last existing header cell) onto the
# new H1 header cell, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" column values for the data rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
